$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add project 8 (Kickstarter Success Classifier) as new row 9.
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Kickstarter Success Classifier "
$ws.Range("C9").Value = "Preditcts if the kickstarter project will be successful "

# Update the skills text on row 5 (Airbnb project): swap "Supervised Machine
# Learning" for "Scikit-Learn" at the end of the skills list.
$ws.Range("D5").Value = "HTML, CSS, Python, PlotlyDash, Flask, plotly, Scikit-Learn"

$ws.Range("D9").Value = "HTML, CSS, Python, Scikit-Learn, Tensorflow, Keras, Flask, Heroku"
$ws.Range("E9").Value = "project8"
$ws.Range("F9").Value = "https://github.com/Kickstarter-Success-Classifier/kickstarter-"
$ws.Range("G9").Value = "https://will-it-kickstart.herokuapp.com/"
$ws.Range("H9").Value = "balle 8 "

$ws.Hyperlinks.Add($ws.Range("F9"), "https://github.com/Kickstarter-Success-Classifier/kickstarter-")
$ws.Hyperlinks.Add($ws.Range("G9"), "https://will-it-kickstart.herokuapp.com/")

$ws.Range("F9").Style = "Hyperlink"
$ws.Range("G9").Style = "Hyperlink"

$ws.Range("J10").Select()
